$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (columns B, E:AB) between rows 49 and 50 ---
$ws.Range("B49").Value = 6865310
$ws.Range("B50").Value = 6865311
$arrR4950A = New-Object "object[,]" 1,24
$arrR4950A[0,0] = "NK Igman Konjic"
$arrR4950A[0,1] = "Zrinjski Mostar"
$arrR4950A[0,2] = 0
$arrR4950A[0,3] = 2
$arrR4950A[0,4] = "A"
$arrR4950A[0,5] = 3.4
$arrR4950A[0,6] = 3.6
$arrR4950A[0,7] = 1.833
$arrR4950A[0,8] = 4.75
$arrR4950A[0,9] = 4.75
$arrR4950A[0,10] = 1.45
$arrR4950A[0,11] = 1.25
$arrR4950A[0,12] = 1.775
$arrR4950A[0,13] = 2.025
$arrR4950A[0,14] = 2.75
$arrR4950A[0,15] = 1.85
$arrR4950A[0,16] = 1.95
$arrR4950A[0,17] = -1
$arrR4950A[0,18] = -1
$arrR4950A[0,19] = 0.45
$arrR4950A[0,20] = -1
$arrR4950A[0,21] = 1.025
$arrR4950A[0,22] = -1
$arrR4950A[0,23] = 0.95
$arrR4950B = New-Object "object[,]" 1,24
$arrR4950B[0,0] = "Sloga"
$arrR4950B[0,1] = "GOSK Gabela"
$arrR4950B[0,2] = 3
$arrR4950B[0,3] = 2
$arrR4950B[0,4] = "H"
$arrR4950B[0,5] = 1.833
$arrR4950B[0,6] = 3.6
$arrR4950B[0,7] = 3.4
$arrR4950B[0,8] = 1.909
$arrR4950B[0,9] = 3.4
$arrR4950B[0,10] = 3.3
$arrR4950B[0,11] = -0.5
$arrR4950B[0,12] = 1.925
$arrR4950B[0,13] = 1.875
$arrR4950B[0,14] = 2.25
$arrR4950B[0,15] = 1.825
$arrR4950B[0,16] = 1.975
$arrR4950B[0,17] = 0.909
$arrR4950B[0,18] = -1
$arrR4950B[0,19] = -1
$arrR4950B[0,20] = 0.925
$arrR4950B[0,21] = -1
$arrR4950B[0,22] = 0.825
$arrR4950B[0,23] = -1
$ws.Range("E49:AB49").Value = $arrR4950A
$ws.Range("E50:AB50").Value = $arrR4950B

# --- Swap match data (columns B, E:AB) between rows 99 and 100 ---
$ws.Range("B99").Value = 6864639
$ws.Range("B100").Value = 6865343
$arrR99100A = New-Object "object[,]" 1,24
$arrR99100A[0,0] = "Zvijezda 09"
$arrR99100A[0,1] = "Borac Banja Luka"
$arrR99100A[0,2] = 1
$arrR99100A[0,3] = 2
$arrR99100A[0,4] = "A"
$arrR99100A[0,5] = 11
$arrR99100A[0,6] = 6
$arrR99100A[0,7] = 1.2
$arrR99100A[0,8] = 10
$arrR99100A[0,9] = 6.5
$arrR99100A[0,10] = 1.181
$arrR99100A[0,11] = 2
$arrR99100A[0,12] = 1.825
$arrR99100A[0,13] = 1.975
$arrR99100A[0,14] = 3
$arrR99100A[0,15] = 1.9
$arrR99100A[0,16] = 1.9
$arrR99100A[0,17] = -1
$arrR99100A[0,18] = -1
$arrR99100A[0,19] = 0.181
$arrR99100A[0,20] = 0.825
$arrR99100A[0,21] = -1
$arrR99100A[0,22] = 0
$arrR99100A[0,23] = 0
$arrR99100B = New-Object "object[,]" 1,24
$arrR99100B[0,0] = "Sloga"
$arrR99100B[0,1] = "NK Posusje"
$arrR99100B[0,2] = 1
$arrR99100B[0,3] = 0
$arrR99100B[0,4] = "H"
$arrR99100B[0,5] = 1.909
$arrR99100B[0,6] = 3.3
$arrR99100B[0,7] = 3.5
$arrR99100B[0,8] = 2.2
$arrR99100B[0,9] = 2.8
$arrR99100B[0,10] = 3.3
$arrR99100B[0,11] = -0.25
$arrR99100B[0,12] = 1.95
$arrR99100B[0,13] = 1.85
$arrR99100B[0,14] = 1.75
$arrR99100B[0,15] = 1.875
$arrR99100B[0,16] = 1.925
$arrR99100B[0,17] = 1.2
$arrR99100B[0,18] = -1
$arrR99100B[0,19] = -1
$arrR99100B[0,20] = 0.95
$arrR99100B[0,21] = -1
$arrR99100B[0,22] = -1
$arrR99100B[0,23] = 0.925
$ws.Range("E99:AB99").Value = $arrR99100A
$ws.Range("E100:AB100").Value = $arrR99100B

# --- Append new rows 183:188 (copy style from row 182, then set values) ---
$ws.Range("A182:AB182").Copy()
$ws.Range("A183:AB188").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rowArr183 = New-Object "object[,]" 1,28
$rowArr183[0,0] = 181
$rowArr183[0,1] = 7952770
$rowArr183[0,2] = "Bosnia Herzegovina Premier Liga"
$rowArr183[0,3] = 45422.625
$rowArr183[0,4] = "FK Sarajevo"
$rowArr183[0,5] = "Zvijezda 09"
$rowArr183[0,6] = 2
$rowArr183[0,7] = 1
$rowArr183[0,8] = "H"
$rowArr183[0,9] = 1.125
$rowArr183[0,10] = 8
$rowArr183[0,11] = 16
$rowArr183[0,12] = 1.142
$rowArr183[0,13] = 9.5
$rowArr183[0,14] = 11
$rowArr183[0,15] = -2.5
$rowArr183[0,16] = 1.925
$rowArr183[0,17] = 1.875
$rowArr183[0,18] = 4
$rowArr183[0,19] = 1.825
$rowArr183[0,20] = 1.975
$rowArr183[0,21] = 0.1419999999999999
$rowArr183[0,22] = -1
$rowArr183[0,23] = -1
$rowArr183[0,24] = -1
$rowArr183[0,25] = 0.875
$rowArr183[0,26] = -1
$rowArr183[0,27] = 0.9750000000000001
$ws.Range("A183:AB183").Value = $rowArr183

$rowArr184 = New-Object "object[,]" 1,28
$rowArr184[0,0] = 182
$rowArr184[0,1] = 7952772
$rowArr184[0,2] = "Bosnia Herzegovina Premier Liga"
$rowArr184[0,3] = 45423.45833333334
$rowArr184[0,4] = "Siroki Brijeg"
$rowArr184[0,5] = "Velez Mostar"
$rowArr184[0,6] = 0
$rowArr184[0,7] = 1
$rowArr184[0,8] = "A"
$rowArr184[0,9] = 3
$rowArr184[0,10] = 3.25
$rowArr184[0,11] = 2.1
$rowArr184[0,12] = 2.9
$rowArr184[0,13] = 3.25
$rowArr184[0,14] = 2.15
$rowArr184[0,15] = 0.25
$rowArr184[0,16] = 1.825
$rowArr184[0,17] = 1.975
$rowArr184[0,18] = 2.25
$rowArr184[0,19] = 2
$rowArr184[0,20] = 1.8
$rowArr184[0,21] = -1
$rowArr184[0,22] = -1
$rowArr184[0,23] = 1.15
$rowArr184[0,24] = -1
$rowArr184[0,25] = 0.9750000000000001
$rowArr184[0,26] = -1
$rowArr184[0,27] = 0.8
$ws.Range("A184:AB184").Value = $rowArr184

$rowArr185 = New-Object "object[,]" 1,28
$rowArr185[0,0] = 183
$rowArr185[0,1] = 7952771
$rowArr185[0,2] = "Bosnia Herzegovina Premier Liga"
$rowArr185[0,3] = 45424.4375
$rowArr185[0,4] = "NK Igman Konjic"
$rowArr185[0,5] = "Zeljeznicar"
$rowArr185[0,6] = 0
$rowArr185[0,7] = 0
$rowArr185[0,8] = "D"
$rowArr185[0,9] = 2.6
$rowArr185[0,10] = 3.1
$rowArr185[0,11] = 2.45
$rowArr185[0,12] = 1.727
$rowArr185[0,13] = 3.6
$rowArr185[0,14] = 4
$rowArr185[0,15] = -0.75
$rowArr185[0,16] = 1.975
$rowArr185[0,17] = 1.825
$rowArr185[0,18] = 2.5
$rowArr185[0,19] = 1.975
$rowArr185[0,20] = 1.825
$rowArr185[0,21] = -1
$rowArr185[0,22] = 2.6
$rowArr185[0,23] = -1
$rowArr185[0,24] = -1
$rowArr185[0,25] = 0.825
$rowArr185[0,26] = -1
$rowArr185[0,27] = 0.825
$ws.Range("A185:AB185").Value = $rowArr185

$rowArr186 = New-Object "object[,]" 1,28
$rowArr186[0,0] = 184
$rowArr186[0,1] = 7952773
$rowArr186[0,2] = "Bosnia Herzegovina Premier Liga"
$rowArr186[0,3] = 45424.60416666666
$rowArr186[0,4] = "NK Posusje"
$rowArr186[0,5] = "Sloga"
$rowArr186[0,6] = 2
$rowArr186[0,7] = 1
$rowArr186[0,8] = "H"
$rowArr186[0,9] = 1.666
$rowArr186[0,10] = 3.25
$rowArr186[0,11] = 4.75
$rowArr186[0,12] = 1.363
$rowArr186[0,13] = 4
$rowArr186[0,14] = 7
$rowArr186[0,15] = -1.25
$rowArr186[0,16] = 1.95
$rowArr186[0,17] = 1.85
$rowArr186[0,18] = 2.25
$rowArr186[0,19] = 1.825
$rowArr186[0,20] = 1.975
$rowArr186[0,21] = 0.363
$rowArr186[0,22] = -1
$rowArr186[0,23] = -1
$rowArr186[0,24] = -0.5
$rowArr186[0,25] = 0.425
$rowArr186[0,26] = 0.825
$rowArr186[0,27] = -1
$ws.Range("A186:AB186").Value = $rowArr186

$rowArr187 = New-Object "object[,]" 1,28
$rowArr187[0,0] = 185
$rowArr187[0,1] = 7952774
$rowArr187[0,2] = "Bosnia Herzegovina Premier Liga"
$rowArr187[0,3] = 45425.5
$rowArr187[0,4] = "FK Tuzla City"
$rowArr187[0,5] = "Borac Banja Luka"
$rowArr187[0,6] = 2
$rowArr187[0,7] = 6
$rowArr187[0,8] = "A"
$rowArr187[0,9] = 5
$rowArr187[0,10] = 4
$rowArr187[0,11] = 1.5
$rowArr187[0,12] = 13
$rowArr187[0,13] = 5.75
$rowArr187[0,14] = 1.166
$rowArr187[0,15] = 2
$rowArr187[0,16] = 1.8
$rowArr187[0,17] = 2
$rowArr187[0,18] = 3.25
$rowArr187[0,19] = 2.025
$rowArr187[0,20] = 1.775
$rowArr187[0,21] = -1
$rowArr187[0,22] = -1
$rowArr187[0,23] = 0.1659999999999999
$rowArr187[0,24] = -1
$rowArr187[0,25] = 1
$rowArr187[0,26] = 1.025
$rowArr187[0,27] = -1
$ws.Range("A187:AB187").Value = $rowArr187

$rowArr188 = New-Object "object[,]" 1,28
$rowArr188[0,0] = 186
$rowArr188[0,1] = 7952775
$rowArr188[0,2] = "Bosnia Herzegovina Premier Liga"
$rowArr188[0,3] = 45425.66666666666
$rowArr188[0,4] = "GOSK Gabela"
$rowArr188[0,5] = "Zrinjski Mostar"
$rowArr188[0,6] = 1
$rowArr188[0,7] = 4
$rowArr188[0,8] = "A"
$rowArr188[0,9] = 6.5
$rowArr188[0,10] = 4.333
$rowArr188[0,11] = 1.363
$rowArr188[0,12] = 5.75
$rowArr188[0,13] = 4.75
$rowArr188[0,14] = 1.38
$rowArr188[0,15] = 1.25
$rowArr188[0,16] = 1.925
$rowArr188[0,17] = 1.875
$rowArr188[0,18] = 2.75
$rowArr188[0,19] = 1.975
$rowArr188[0,20] = 1.825
$rowArr188[0,21] = -1
$rowArr188[0,22] = -1
$rowArr188[0,23] = 0.3799999999999999
$rowArr188[0,24] = -1
$rowArr188[0,25] = 0.875
$rowArr188[0,26] = 0.9750000000000001
$rowArr188[0,27] = -1
$ws.Range("A188:AB188").Value = $rowArr188

